$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.446.51'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '1.826.91'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'316.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("D7").Value = "'0.5152"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.68%  '
$ws.Range("D8").Value = "'0.3932"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.84%  '
$ws.Range("D9").Value = "'0.07711"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.78%  '
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").Value = "'1.114"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").Value = "'21.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.71%  '
$ws.Range("D13").Value = "'6.279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").Value = "'7.575"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.41%  '
$ws.Range("B15").Value = 'BinanceUSD'
$ws.Range("C15").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D15").Value = "'1.002"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.35%  '
$ws.Range("D16").Value = '1.824.37'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = "'93.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.14%  '
$ws.Range("D18").Value = "'0.00001080"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("D19").Value = "'0.06624"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.97%  '
$ws.Range("D20").Value = "'17.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.42%  '
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = "'6.064"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("D23").Value = '28.459.24'
$ws.Range("E23").Value = '  +1.84%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = "'2.242"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.09%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = "'2.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.89%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = "'157.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").Value = "'20.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.47%  '
$ws.Range("D29").Value = '2.035.24'
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("D30").Value = "'124.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.29%  '
$ws.Range("D31").Value = "'1.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.24%  '
$ws.Range("D32").Value = "'0.1098"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.52%  '
$ws.Range("D33").Value = "'5.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.31%  '
$ws.Range("D34").Value = "'3.614"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("D35").Value = "'0.07182"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.64%  '
$ws.Range("D36").Value = "'0.2233"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("D37").Value = "'8.982"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.94%  '
$ws.Range("D38").Value = "'0.02321"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.82%  '
$ws.Range("D39").Value = "'5.149"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6244"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.09%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = "'11.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("E42").Value = '  +2.24%  '
$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").Value = "'1.395"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.61%  '
$ws.Range("D45").Value = "'13.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").Value = "'0.5901"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.56%  '
$ws.Range("D47").Value = "'3.708"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("D48").Value = "'124.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("D49").Value = "'1.979"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.66%  '
$ws.Range("D50").Value = "'1.181"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.52%  '
$ws.Range("D51").Value = "'0.06927"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.92%  '
